$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 378, shifting rows 378:413 down to 379:414
$ws.Rows.Item(378).Insert()

# Populate the newly inserted row 378 with the new record's data
$ws.Cells.Item(378, 1).Value = 5
$ws.Cells.Item(378, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(378, 3).Value = "Maule"
$ws.Cells.Item(378, 4).Value = 45106
$ws.Cells.Item(378, 5).Value = 7
$ws.Cells.Item(378, 6).Value = "Fruta"
$ws.Cells.Item(378, 7).Value = 100108
$ws.Cells.Item(378, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(378, 9).Value = 100108005
$ws.Cells.Item(378, 10).Value = "Piña"
$ws.Cells.Item(378, 11).Value = "Caramelo"
$ws.Cells.Item(378, 12).Value = "Segunda"
$ws.Cells.Item(378, 13).Value = 150
$ws.Cells.Item(378, 14).Value = 24000
$ws.Cells.Item(378, 15).Value = 24000
$ws.Cells.Item(378, 16).Value = 24000
$ws.Cells.Item(378, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(378, 18).Value = "Ecuador"
$ws.Cells.Item(378, 19).Value = 1714
$ws.Cells.Item(378, 20).Value = 14
